$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Remove existing threaded comments (on B4 and A11) - this also drops the
# legacy comments/threaded comments parts and the legacyDrawing ref.
$c1 = $ws2.Range("B4").Comment
if ($c1 -ne $null) { $c1.Delete() }
$c2 = $ws2.Range("A11").Comment
if ($c2 -ne $null) { $c2.Delete() }

# Column widths (col B wider, new col C added)
$ws2.Columns.Item(2).ColumnWidth = 25
$ws2.Columns.Item(3).ColumnWidth = 53.333333333333336

# --- Rewrite the parameter table rows 4-13 (rows 4-12 shuffled/extended) ---
# Labels / values that reuse already-existing shared strings - order doesn't matter.

# Row 4: CurrentYear (was Look Ahead w/ formula)
$ws2.Range("A4").Value = "CurrentYear"
$ws2.Range("B4").Value = 0

# Row 5: InvestmentIteration (new position)
$ws2.Range("A5").Value = "InvestmentIteration"
$ws2.Range("B5").Value = 0

# Row 6: pastTimeHorizon (wrap-text style)
$ws2.Range("A6").Value = "pastTimeHorizon"
$ws2.Range("A6").WrapText = $true
$ws2.Range("B6").Value = 3

# Row 7: Look Ahead (now points to B13)
$ws2.Range("A7").Value = "Look Ahead"
$ws2.Range("B7").Formula = "=B13"

# Row 8: Country
$ws2.Range("A8").Value = "Country"
$ws2.Range("B8").Value = "DE"

# Row 9: short_term_investment_minimal_irr
$ws2.Range("A9").Value = "short_term_investment_minimal_irr"
$ws2.Range("B9").Value = 0.3

# Row 10: start_year_fuel_trends
$ws2.Range("A10").Value = "start_year_fuel_trends"
$ws2.Range("B10").Value = 5

# Row 12: maximum_investment_capacity_per_year
$ws2.Range("A12").Style = "Normal"
$ws2.Range("A12").Value = "maximum_investment_capacity_per_year"
$ws2.Range("B12").Value = 10000

# --- New strings, written in the exact order they were first introduced ---

# Row 11: start_year_dismantling (brand new row/label)
$ws2.Range("A11").Value = "start_year_dismantling"
$ws2.Range("B11").Value = 4

# Row 13: max_permit_build_time (moved here) + note
$ws2.Range("A13").Value = "max_permit_build_time"
$ws2.Range("B13").Value = 4
$ws2.Range("C13").Value = "should be maximum permit and lead time of candidate technologies(check emlab parameters) "

# Row 7 note
$ws2.Range("C7").Value = "Be sure that there is data ready until 'End Year' + 'Look Ahead'"

# Row 9 note
$ws2.Range("C9").Value = "Minimal IRR to make investment decisions"

# Row 12 note
$ws2.Range("C12").Value = "MW"

# Row 10 note
$ws2.Range("C10").Value = "Year when the prices are not longer interpolated, but determined through trend"

# Row 11 note
$ws2.Range("C11").Value = "Year when the dismantling begins, based on the profits"

# Row 5 note
$ws2.Range("C5").Value = "Count of number of iterations per year. This is changed by the program"

# Row 6 note + validation formula
$ws2.Range("C6").Value = "to calculate the profits for dismantling. Should be smaller than start year dismantling"
$ws2.Range("D6").Formula = '=IF(B6>B11," !!! past time horizon should be at least the year of dismantling","ok")'

# Conditional formatting on D6: highlight red when not "ok"
$fc = $ws2.Range("D6").FormatConditions.Add([Microsoft.Office.Interop.Excel.XlFormatConditionType]::xlCellValue, [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlNotEqual, '"ok"')
$fc.Interior.Color = 255

# Update selection to match saved view state
$ws2.Range("C15").Select()

$wb.Save()
